$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 195
$ws.Range("A195").Formula = '''6/7/2019'
$ws.Range("B195").Value = 'SPE7M5-19-V-9760'
$ws.Range("C195").Formula = '''13'
$ws.Range("D195").Formula = '''$1,608.36'
$ws.Range("E195").Formula = '''5935016151446'
$ws.Range("F195").Value = 'BACKSHELL,ELECTRICA'
$ws.Range("G195").Value = 'Glenair'
$ws.Range("H195").Value = '447HS325XW1106B'
$ws.Range("I195").Value = 'CP'
$ws.Range("J195").Value = '2019 NOV 25'

# Row 196
$ws.Range("A196").Formula = '''6/7/2019'
$ws.Range("B196").Value = 'SPE7L7-19-V-1150'
$ws.Range("C196").Formula = '''1'
$ws.Range("D196").Formula = '''$2,394.46'
$ws.Range("E196").Formula = '''6130014355672'
$ws.Range("F196").Value = 'POWER SUPPLY'
$ws.Range("G196").Value = 'Druck'
$ws.Range("H196").Value = 'ADTS405-1892-25-M0'
$ws.Range("I196").Value = 'M41'
$ws.Range("J196").Value = '2019 NOV 14'

